$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Globo", "Bom Dia Brasil", "Cultura", "2025-03-28T15:57", "Neutro", "Sem Nota", "aaa"),
    @("Globo", "Bom Dia Rio", "Codemca", "2025-03-28T15:57", "Neutro", "Sem Nota", "teste2"),
    @("Record", "Balanço Geral", "Codemca", "2025-03-14T16:08", "Neutro", "Com Nota", "a1"),
    @("Globo", "Inter TV Rural", "Agricultura", "2025-03-28T18:11", "Neutro", "Sem Nota", "aaaaaaaaaaaaaaaaaaaaaaateste3333333"),
    @("Globo", "Inter TV Rural", "Cultura", "2025-03-28T16:12", "Neutro", "Com Nota", "212121212")
)

$startRow = 8

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($col = 1; $col -le 7; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $value = $rowData[$col - 1]

        # Guard against Excel auto-converting purely numeric-looking text
        # (e.g. "212121212") into a real number. Force the cell to Text
        # format first so the value is stored verbatim as a string.
        if ($value -match '^-?\d+(\.\d+)?$') {
            $cell.NumberFormat = "@"
        }

        $cell.Value = $value
    }
}
